$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: "HS Percentile" ---
# Copy the header style from I1 (bold / bordered / centered) onto J1, then set its text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "HS Percentile"

# Data values for the new column
$ws.Range("J2").Value = 96
$ws.Range("J3").Value = 31

# --- Row 3 edits ---
# A3 and C3 become the text "1" (quote-prefixed / text-formatted), but keep the
# default (unstyled) cell style, matching the target workbook.
$a3 = $ws.Range("A3")
$a3.NumberFormat = "@"
$a3.Value = "1"
$a3.Style = "Normal"

$c3 = $ws.Range("C3")
$c3.NumberFormat = "@"
$c3.Value = "1"
$c3.Style = "Normal"

# Numeric columns E3:I3 all reset to 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
